# Apply updated crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.048.27"
$ws.Range("E2").Value = "'  +3.87%  "
$ws.Range("D3").Value = "'3.212.70"
$ws.Range("E3").Value = "'  +1.18%  "
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'219.18"
$ws.Range("E5").Value = "'  +5.96%  "
$ws.Range("D6").Value = "'647.51"
$ws.Range("E6").Value = "'  +5.84%  "
$ws.Range("E7").Value = "'  +4.11%  "
$ws.Range("D8").Value = "'0.704"
$ws.Range("E8").Value = "'  +4.70%  "
$ws.Range("E9").Value = "'  +0.06%  "
$ws.Range("D10").Value = "'3.209.60"
$ws.Range("E10").Value = "'  +1.18%  "
$ws.Range("D11").Value = "'0.580"
$ws.Range("E11").Value = "'  +7.86%  "
$ws.Range("D12").Value = "'0.182"
$ws.Range("E12").Value = "'  +2.24%  "
$ws.Range("E13").Value = "'  +6.10%  "
$ws.Range("D14").Value = "'5.44"
$ws.Range("E14").Value = "'  +3.21%  "
$ws.Range("D15").Value = "'33.61"
$ws.Range("E15").Value = "'  +4.39%  "
$ws.Range("D16").Value = "'90.745.01"
$ws.Range("E16").Value = "'  +3.89%  "
$ws.Range("D17").Value = "'3.801.03"
$ws.Range("E17").Value = "'  +1.31%  "
$ws.Range("D18").Value = "'3.220.79"
$ws.Range("E18").Value = "'  +2.42%  "
$ws.Range("B19").Value = "'SuiNetwork"
$ws.Range("C19").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").Value = "'3.39"
$ws.Range("E19").Value = "'  +12.28%  "
$ws.Range("B20").Value = "'PEPE"
$ws.Range("C20").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D20").Value = "'0.0000226"
$ws.Range("E20").Value = "'  +72.39%  "
$ws.Range("D21").Value = "'13.58"
$ws.Range("E21").Value = "'  +0.76%  "
$ws.Range("D22").Value = "'442.11"
$ws.Range("E22").Value = "'  +6.00%  "
$ws.Range("D23").Value = "'8.72"
$ws.Range("E23").Value = "'  +2.32%  "
$ws.Range("D24").Value = "'5.12"
$ws.Range("E24").Value = "'  +0.24%  "
$ws.Range("D25").Value = "'5.32"
$ws.Range("E25").Value = "'  +2.16%  "
$ws.Range("D26").Value = "'11.92"
$ws.Range("E26").Value = "'  -0.13%  "
$ws.Range("D27").Value = "'82.06"
$ws.Range("E27").Value = "'  +11.40%  "
$ws.Range("D28").Value = "'3.383.63"
$ws.Range("E28").Value = "'  +1.55%  "
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("E30").Value = "'  +0.42%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  +0.24%  "
$ws.Range("D32").Value = "'4.19"
$ws.Range("E32").Value = "'  +39.60%  "
$ws.Range("B33").Value = "'InternetComputer(DFINITY)"
$ws.Range("C33").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'8.52"
$ws.Range("E33").Value = "'  +2.86%  "
$ws.Range("B34").Value = "'Bittensor"
$ws.Range("C34").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'546.26"
$ws.Range("E34").Value = "'  +0.08%  "
$ws.Range("D35").Value = "'7.16"
$ws.Range("E35").Value = "'  +6.18%  "
$ws.Range("D36").Value = "'1.94"
$ws.Range("E36").Value = "'  +5.01%  "
$ws.Range("E37").Value = "'  +0.38%  "
$ws.Range("D38").Value = "'22.61"
$ws.Range("E38").Value = "'  +3.26%  "
$ws.Range("D39").Value = "'22.42"
$ws.Range("E39").Value = "'  +2.79%  "
$ws.Range("D40").Value = "'0.127"
$ws.Range("E40").Value = "'  -3.71%  "
$ws.Range("E41").Value = "'  +0.14%  "
$ws.Range("D42").Value = "'1.95"
$ws.Range("E42").Value = "'  +2.93%  "
$ws.Range("E43").Value = "'  -0.03%  "
$ws.Range("D44").Value = "'0.377"
$ws.Range("E44").Value = "'  +1.39%  "
$ws.Range("D45").Value = "'45.19"
$ws.Range("E45").Value = "'  +4.38%  "
$ws.Range("D46").Value = "'146.50"
$ws.Range("E46").Value = "'  -0.51%  "
$ws.Range("D47").Value = "'174.43"
$ws.Range("E47").Value = "'  +0.28%  "
$ws.Range("D48").Value = "'0.766"
$ws.Range("E48").Value = "'  +9.30%  "
$ws.Range("E49").Value = "'  -1.13%  "
$ws.Range("D50").Value = "'1.25"
$ws.Range("E50").Value = "'  +0.78%  "
$ws.Range("E51").Value = "'  +6.92%  "
